$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "ImageUrl"
$ws.Range("M2").Value = "TomatoCurry.jpg"
$ws.Range("M3").Value = "Brinjal-rice.jpg"
$ws.Range("M4").Value = "gajar.jpg"

$ws.Range("M2:M4").WrapText = $true

$ws.Range("M4").Select()
